$d = $word.ActiveDocument

$p3 = $d.Paragraphs(3).Range
$p3.Find.Execute("System should be able to flag patients in need of a dietitian using their data.", $true, $false, $false, $false, $false,
                  $true, 1, $false, "Create accounts and delete (staff)", 2)

$p4 = $d.Paragraphs(4).Range
$p4.Find.Execute("Create accounts and delete (staff)", $true, $false, $false, $false, $false,
                  $true, 1, $false, "It should be able to run on different platforms (Windows, Mac, Linux)", 2)
